# Re-Run Stats + Tidy Up DiCE/SHAP
# Adds a "Mean" summary row (row 25) under the LIME XAI metrics table,
# averaging each metric column over the 20 sample rows (3:22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Label cell for the new summary row
$ws.Range("B25").Value = "Mean"

# Column C gets its own (non-shared) AVERAGE formula; D:G share the same
# formula pattern (relative reference), matching how Excel's fill/copy
# creates a shared formula group when you drag the first one across.
$ws.Range("C25").Formula = "=AVERAGE(C3:C22)"
$ws.Range("D25:G25").Formula = "=AVERAGE(D3:D22)"

# Leave the selection where the author's Excel session ended up.
$ws.Range("D31").Select() | Out-Null
